$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Locate the paragraph whose entire content is a single tab character
# (it follows the "git push -u origin master" list item and is immediately
# followed by the document's final, empty paragraph).
$tabParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`t`r") {
        $tabParaIndex = $i
        break
    }
}

$tabPara = $d.Paragraphs.Item($tabParaIndex)

# 1) Insert a brand-new, empty paragraph styled "ListParagraph" (with a
#    1440-twip left indent) right before the tab-only paragraph.
$insertionPoint = $d.Range($tabPara.Range.Start, $tabPara.Range.Start)
$insertionPoint.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:ind w:left='1440'/><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr></w:p>") | Out-Null

# The tab-only paragraph shifted down by one position.
$tabParaIndex = $tabParaIndex + 1
$tabPara = $d.Paragraphs.Item($tabParaIndex)

# 2) Give that tab-only paragraph a second run containing "clone a folder"
#    (keeping the original tab run intact as its own run).
$tabPara.Range.InsertXML("<w:p $wns><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>clone a folder</w:t></w:r></w:p>") | Out-Null

# 3) Turn the document's final (empty) paragraph into a new paragraph
#    containing a tab followed by "git clone".
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertXML("<w:p $wns><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:tab/><w:t>git clone</w:t></w:r></w:p>") | Out-Null
